# edit.ps1 - applies the "15/20: edit another user (as admin)" revision
# to Resources/Funcionalidades.docx
#
# Changes:
#  1. Header date "14/19" -> "15/20"
#  2. "Editar usuario AU" cell split into "Editar " / "mi propio perfil" / " AU"
#  3. New table row "Editar otro usuario A" inserted right before the
#     "Borrar usuario A" row
#  4. The "Esto tiene que ser como ADMIN..." remark removed from the
#     "Borrar usuario A" row (4th cell becomes an empty paragraph)

$d = $word.ActiveDocument

# Reusable wrapper: build a WordProcessingML "single file package" snippet
# that InsertXML understands, given a body fragment.
function New-PkgXml([string]$bodyFragment) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Header date: "14/19" -> "15/20" (runs: "1","4","/1","9")
# ---------------------------------------------------------------------
$datePara = $d.Paragraphs(3)
$dateStart = $datePara.Range.Start

$d.Range($dateStart + 1, $dateStart + 2).Text = "5"          # "4"  -> "5"
$d.Range($dateStart + 2, $dateStart + 4).Text = "/"           # "/1" -> "/"
$d.Range($dateStart + 3, $dateStart + 4).Text = "20"          # "9"  -> "20"

# ---------------------------------------------------------------------
# 2) "Editar usuario AU" -> "Editar " + "mi propio perfil" + " AU"
#    (kept as three separate highlighted runs)
# ---------------------------------------------------------------------
$editRange = $d.Content
$found = $editRange.Find.Execute("Editar usuario AU", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitXml = '<w:p>' +
        '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve">Editar </w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>mi propio perfil</w:t></w:r>' +
        '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> AU</w:t></w:r>' +
        '</w:p>'
    $editRange.InsertXML((New-PkgXml $splitXml))
}

# ---------------------------------------------------------------------
# 3) Insert a new row "Editar otro usuario A" right before the
#    "Borrar usuario A" row (i.e. right after "Ver todos los usuarios A")
# ---------------------------------------------------------------------
$table = $d.Tables(1)
$targetRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $rowText = $table.Rows.Item($i).Range.Text
    if ($rowText -like "Borrar usuario A*") {
        $targetRow = $table.Rows.Item($i)
        break
    }
}

$newRow = $table.Rows.Add($targetRow)

$firstCellXml = '<w:p><w:pPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Editar otro usuario A</w:t></w:r></w:p>'
$newRow.Cells.Item(1).Range.InsertXML((New-PkgXml $firstCellXml))

# ---------------------------------------------------------------------
# 4) Remove the "Esto tiene que ser como ADMIN..." remark text
# ---------------------------------------------------------------------
$remarkRange = $d.Content
$foundRemark = $remarkRange.Find.Execute("Esto tiene que ser como ADMIN con la lista de usuarios por delante", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundRemark) {
    $remarkRange.Text = ""
}
